# Coin list refresh: update Price (D) and Volume(1h) (E) text columns.
# Cells whose new text parses as a plain number (single decimal point, e.g.
# "240.80") need a leading apostrophe so Excel stores them as text rather
# than silently converting them to a numeric value -- exactly like typing
# '240.80 into the cell by hand. Cells with dotted thousands separators
# (e.g. "29.429.96") are never parsed as numbers, so no apostrophe is needed.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.429.96'
$ws.Range("E2").Value = '  +0.34%  '

$ws.Range("D3").Value = '1.849.03'
$ws.Range("E3").Value = '  +0.36%  '

$ws.Range("D5").Value = '''240.80'
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").Value = '''0.6297'
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").Value = '''1.001'
$ws.Range("E7").Value = '  +0.09%  '

$ws.Range("D8").Value = '''0.07698'
$ws.Range("E8").Value = '  +2.33%  '

$ws.Range("D9").Value = '''0.2923'
$ws.Range("E9").Value = '  -0.52%  '

$ws.Range("D10").Value = '''24.75'
$ws.Range("E10").Value = '  +1.31%  '

$ws.Range("D11").Value = '''0.07743'
$ws.Range("E11").Value = '  +0.72%  '

$ws.Range("D12").Value = '1.847.38'
$ws.Range("E12").Value = '  +0.64%  '

$ws.Range("D13").Value = '''5.030'
$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("D14").Value = '''0.6799'
$ws.Range("E14").Value = '  +0.33%  '

$ws.Range("D15").Value = '''0.00001070'
$ws.Range("E15").Value = '  +1.62%  '

$ws.Range("D16").Value = '''83.59'
$ws.Range("E16").Value = '  +0.83%  '

$ws.Range("D17").Value = '''6.179'
$ws.Range("E17").Value = '  +1.06%  '

$ws.Range("D18").Value = '29.446.88'
$ws.Range("E18").Value = '  +0.28%  '

$ws.Range("D19").Value = '''228.11'
$ws.Range("E19").Value = '  +0.05%  '

$ws.Range("D20").Value = '''12.43'
$ws.Range("E20").Value = '  +0.08%  '

$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").Value = '''7.426'
$ws.Range("E22").Value = '  +0.13%  '

$ws.Range("E23").Value = '  +0.10%  '

$ws.Range("D24").Value = '''157.61'
$ws.Range("E24").Value = '  +0.60%  '

$ws.Range("D25").Value = '''0.1381'
$ws.Range("E25").Value = '  -0.40%  '

$ws.Range("D26").Value = '''8.402'
$ws.Range("E26").Value = '  +0.65%  '

$ws.Range("E27").Value = '  +0.57%  '

$ws.Range("D28").Value = '''1.354'
$ws.Range("E28").Value = '  +6.15%  '

$ws.Range("D29").Value = '''1.467'
$ws.Range("E29").Value = '  +1.05%  '

$ws.Range("D30").Value = '''0.05679'
$ws.Range("E30").Value = '  +1.05%  '

$ws.Range("D31").Value = '''4.119'
$ws.Range("E31").Value = '  +0.48%  '

$ws.Range("D32").Value = '''4.028'
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("E33").Value = '  +0.85%  '

$ws.Range("E34").Value = '  +0.53%  '

$ws.Range("E35").Value = '  -0.11%  '

$ws.Range("E36").Value = '  -0.20%  '

$ws.Range("D37").Value = '''2.777'
$ws.Range("E37").Value = '  +0.54%  '

$ws.Range("E38").Value = '  -1.01%  '

$ws.Range("D39").Value = '1.220.12'
$ws.Range("E39").Value = '  -1.60%  '

$ws.Range("D40").Value = '''6.545'
$ws.Range("E40").Value = '  +4.94%  '

$ws.Range("D41").Value = '''0.9065'
$ws.Range("E41").Value = '  +0.71%  '

$ws.Range("E42").Value = '  +0.17%  '

$ws.Range("D43").Value = '''101.76'
$ws.Range("E43").Value = '  -0.09%  '

$ws.Range("D44").Value = '''66.18'
$ws.Range("E44").Value = '  +1.04%  '

$ws.Range("D45").Value = '''0.00000000120'
$ws.Range("E45").Value = '  +1.09%  '

$ws.Range("D46").Value = '''7.142'
$ws.Range("E46").Value = '  +0.87%  '

$ws.Range("E47").Value = '  +0.83%  '

$ws.Range("D48").Value = '''9.033'
$ws.Range("E48").Value = '  +1.06%  '

$ws.Range("E49").Value = '  +2.45%  '

$ws.Range("D50").Value = '''1.674'
$ws.Range("E50").Value = '  +0.44%  '

$ws.Range("E51").Value = '  +0.12%  '
